$d = $word.ActiveDocument

# Phase 1: replace each original value with a unique placeholder token
# to avoid collisions between an old value and another cell's new value.
$d.Content.Find.Execute("74÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@0@@", 2) | Out-Null
$d.Content.Find.Execute("54÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@1@@", 2) | Out-Null
$d.Content.Find.Execute("92÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@2@@", 2) | Out-Null
$d.Content.Find.Execute("17÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@3@@", 2) | Out-Null
$d.Content.Find.Execute("91÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@4@@", 2) | Out-Null
$d.Content.Find.Execute("95÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@5@@", 2) | Out-Null
$d.Content.Find.Execute("41÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@6@@", 2) | Out-Null
$d.Content.Find.Execute("30÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@7@@", 2) | Out-Null
$d.Content.Find.Execute("62÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@8@@", 2) | Out-Null
$d.Content.Find.Execute("57÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@9@@", 2) | Out-Null
$d.Content.Find.Execute("71÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@10@@", 2) | Out-Null
$d.Content.Find.Execute("26÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@11@@", 2) | Out-Null
$d.Content.Find.Execute("15÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@12@@", 2) | Out-Null
$d.Content.Find.Execute("60÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@13@@", 2) | Out-Null
$d.Content.Find.Execute("35÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@14@@", 2) | Out-Null
$d.Content.Find.Execute("16÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@15@@", 2) | Out-Null
$d.Content.Find.Execute("83÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@16@@", 2) | Out-Null
$d.Content.Find.Execute("70÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@17@@", 2) | Out-Null
$d.Content.Find.Execute("51÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@18@@", 2) | Out-Null
$d.Content.Find.Execute("41÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@19@@", 2) | Out-Null
$d.Content.Find.Execute("64÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@20@@", 2) | Out-Null
$d.Content.Find.Execute("56÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@21@@", 2) | Out-Null
$d.Content.Find.Execute("18÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@22@@", 2) | Out-Null
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@23@@", 2) | Out-Null
$d.Content.Find.Execute("57÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@24@@", 2) | Out-Null

# Phase 2: replace placeholders with final values
$d.Content.Find.Execute("@@0@@", $true, $false, $false, $false, $false, $true, 1, $false, "86÷3=", 2) | Out-Null
$d.Content.Find.Execute("@@1@@", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@2@@", $true, $false, $false, $false, $false, $true, 1, $false, "17÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@3@@", $true, $false, $false, $false, $false, $true, 1, $false, "26÷7=", 2) | Out-Null
$d.Content.Find.Execute("@@4@@", $true, $false, $false, $false, $false, $true, 1, $false, "50÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@5@@", $true, $false, $false, $false, $false, $true, 1, $false, "11÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@6@@", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@7@@", $true, $false, $false, $false, $false, $true, 1, $false, "59÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@8@@", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@9@@", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@10@@", $true, $false, $false, $false, $false, $true, 1, $false, "63÷7=", 2) | Out-Null
$d.Content.Find.Execute("@@11@@", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=", 2) | Out-Null
$d.Content.Find.Execute("@@12@@", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@13@@", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@14@@", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@15@@", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=", 2) | Out-Null
$d.Content.Find.Execute("@@16@@", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@17@@", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=", 2) | Out-Null
$d.Content.Find.Execute("@@18@@", $true, $false, $false, $false, $false, $true, 1, $false, "10÷6=", 2) | Out-Null
$d.Content.Find.Execute("@@19@@", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=", 2) | Out-Null
$d.Content.Find.Execute("@@20@@", $true, $false, $false, $false, $false, $true, 1, $false, "14÷4=", 2) | Out-Null
$d.Content.Find.Execute("@@21@@", $true, $false, $false, $false, $false, $true, 1, $false, "89÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@22@@", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=", 2) | Out-Null
$d.Content.Find.Execute("@@23@@", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=", 2) | Out-Null
$d.Content.Find.Execute("@@24@@", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=", 2) | Out-Null
